{"js": "// Apply the five text edits described by the diff.\n// Each edit is located with a unique, unambiguous search string and then\n// replaced in place with Range.insertText(..., \"Replace\"), which keeps the\n// surrounding run formatting intact.\n\nconst body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) Oliver's paragraph: the visible text does not change, but the source\n//    XML merges the runs that used to wrap \"data\" in proofErr/gramStart\n//    markers into a single run. Re-inserting the identical text over the\n//    full sentence collapses those runs the same way.\nawait replaceOnce(\n  \": My ideal job has not changed after reading the Burning Glass Technologies data. My goal is to be a Game Developer/Software Engineer. After seeing the data I have noticed that a lot of the higher up skill sets are related to skills I believe a software engineer should have and a skills that employers are looking for, leading me to believe following my career choice could create amazing opportunities in the future. Therefore, I am not changing my ideal job and do not believe it is changing any time soon.\",\n  \": My ideal job has not changed after reading the Burning Glass Technologies data. My goal is to be a Game Developer/Software Engineer. After seeing the data I have noticed that a lot of the higher up skill sets are related to skills I believe a software engineer should have and a skills that employers are looking for, leading me to believe following my career choice could create amazing opportunities in the future. Therefore, I am not changing my ideal job and do not believe it is changing any time soon.\"\n);\n\n// 2) Corbin: \"requires earning merit\" -> \"requires earnt merit\"\nawait replaceOnce(\n  \"Such a position requires earning merit through\",\n  \"Such a position requires earnt merit through\"\n);\n\n// 3) Corbin: \"specialised skills in one field\" -> \"specialised skills in only one field\"\nawait replaceOnce(\n  \"then I would only have specialised skills in one field.\",\n  \"then I would only have specialised skills in only one field.\"\n);\n\n// 4) Corbin: \"hopefully open doorways, maybe.\" -> \"hopefully open doorways, possibly.\"\nawait replaceOnce(\n  \"hopefully open doorways, maybe.\",\n  \"hopefully open doorways, possibly.\"\n);\n\n// 5) CTO paragraph: \"be in charge of the organisations technical needs\"\n//    -> \"be in charge of the organisation's technical needs\"\nawait replaceOnce(\n  \"I would be in charge of the organisations technical needs\",\n  \"I would be in charge of the organisation\\u2019s technical needs\"\n);\n", "ps1": "# Apply the five text edits described by the diff using Word's Find/Replace\n# COM automation. Each call targets a unique, unambiguous phrase so the\n# replacement lands exactly where the diff shows it.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $result = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"Find/Replace did not find text: $oldText\"\n    }\n}\n\n# 1) Oliver's paragraph: the visible text does not change, but the source\n#    XML merges the runs that used to wrap \"data\" in proofErr/gramStart\n#    markers into a single run. Re-running Find/Replace with identical\n#    text over the full sentence collapses those runs the same way.\n$oliverText = \": My ideal job has not changed after reading the Burning Glass Technologies data. My goal is to be a Game Developer/Software Engineer. After seeing the data I have noticed that a lot of the higher up skill sets are related to skills I believe a software engineer should have and a skills that employers are looking for, leading me to believe following my career choice could create amazing opportunities in the future. Therefore, I am not changing my ideal job and do not believe it is changing any time soon.\"\nReplace-Text $oliverText $oliverText\n\n# 2) Corbin: \"requires earning merit\" -> \"requires earnt merit\"\nReplace-Text \"Such a position requires earning merit through\" \"Such a position requires earnt merit through\"\n\n# 3) Corbin: \"specialised skills in one field\" -> \"specialised skills in only one field\"\nReplace-Text \"then I would only have specialised skills in one field.\" \"then I would only have specialised skills in only one field.\"\n\n# 4) Corbin: \"hopefully open doorways, maybe.\" -> \"hopefully open doorways, possibly.\"\nReplace-Text \"hopefully open doorways, maybe.\" \"hopefully open doorways, possibly.\"\n\n# 5) CTO paragraph: \"be in charge of the organisations technical needs\"\n#    -> \"be in charge of the organisation's technical needs\"\nReplace-Text \"I would be in charge of the organisations technical needs\" \"I would be in charge of the organisation\u2019s technical needs\"\n"}
